$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 01:52"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 957505
$ws.Range("C4").Value = 32273
$ws.Range("D4").Value = 116201
$ws.Range("E4").Value = 787144
$ws.Range("F4").Value = 15110
$ws.Range("G4").Value = 1967
$ws.Range("H4").Value = 54160

# Row 31: Pakistan
$ws.Range("A31").Value = "Pakistan"
$ws.Range("B31").Value = 12723
$ws.Range("C31").Value = 783
$ws.Range("D31").Value = 2866
$ws.Range("E31").Value = 9588
$ws.Range("F31").Value = 111
$ws.Range("G31").Value = 16
$ws.Range("H31").Value = 269

# Row 32: Singapur
$ws.Range("A32").Value = "Singapur"
$ws.Range("B32").Value = 12693
$ws.Range("C32").Value = 618
$ws.Range("D32").Value = 1002
$ws.Range("E32").Value = 11679
$ws.Range("F32").Value = 24
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 12

# Row 56: Argentina
$ws.Range("A56").Value = "Argentina"
$ws.Range("B56").Value = 3780
$ws.Range("C56").Value = 173
$ws.Range("D56").Value = 1030
$ws.Range("E56").Value = 2565
$ws.Range("F56").Value = 144
$ws.Range("G56").Value = 9
$ws.Range("H56").Value = 185

# Row 57: Luxemburgo
$ws.Range("A57").Value = "Luxemburgo"
$ws.Range("B57").Value = 3711
$ws.Range("C57").Value = 16
$ws.Range("D57").Value = 3088
$ws.Range("E57").Value = 538
$ws.Range("F57").Value = 26
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 85

# Row 120: Venezuela
$ws.Range("A120").Value = "Venezuela"
$ws.Range("B120").Value = 323
$ws.Range("C120").Value = 5
$ws.Range("D120").Value = 132
$ws.Range("E120").Value = 181
$ws.Range("F120").Value = 3
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 10

# Row 121: Montenegro
$ws.Range("A121").Value = "Montenegro"
$ws.Range("B121").Value = 320
$ws.Range("C121").Value = 1
$ws.Range("D121").Value = 153
$ws.Range("E121").Value = 161
$ws.Range("F121").Value = 7
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 6

# Row 133: Maldivas
$ws.Range("A133").Value = "Maldivas"
$ws.Range("B133").Value = 177
$ws.Range("C133").Value = 48
$ws.Range("D133").Value = 17
$ws.Range("E133").Value = 160
$ws.Range("F133").Value = 2
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 0

# Row 134: Gabon
$ws.Range("A134").Value = "Gabon"
$ws.Range("B134").Value = 176
$ws.Range("C134").Value = 4
$ws.Range("D134").Value = 30
$ws.Range("E134").Value = 143
$ws.Range("F134").Value = 1
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 3

# Row 135: Martinica
$ws.Range("A135").Value = "Martinica"
$ws.Range("B135").Value = 175
$ws.Range("C135").Value = 5
$ws.Range("D135").Value = 77
$ws.Range("E135").Value = 84
$ws.Range("F135").Value = 7
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 14

# Row 136: Guadalupe
$ws.Range("A136").Value = "Guadalupe"
$ws.Range("B136").Value = 149
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 82
$ws.Range("E136").Value = 55
$ws.Range("F136").Value = 11
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 12

# Row 137: Birmania
$ws.Range("A137").Value = "Birmania"
$ws.Range("B137").Value = 144
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 10
$ws.Range("E137").Value = 129
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 5
